# Updated symbol list on Tue Dec 13 21:36:50 UTC 2022 with GitHub Actions
#
# Applies price (column D) updates, plus two coin-row swaps (rows 42/43
# and rows 49/50) to the "cryptos" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    # Force the cell to stay a TEXT cell even when the string looks like a
    # number (e.g. "269.29"), matching the workbook's existing convention
    # of storing all Price/Volume columns as inline/shared strings.
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# ----- Column D (Price) numeric-looking text updates -----
Set-TextValue "D2"  "269.29"
Set-TextValue "D3"  "22.92"
Set-TextValue "D4"  "6.329"
Set-TextValue "D5"  "0.06198"
Set-TextValue "D6"  "3.642"
Set-TextValue "D8"  "1.392"
Set-TextValue "D11" "0.1605"
Set-TextValue "D12" "0.08284"
Set-TextValue "D13" "0.03472"
Set-TextValue "D15" "0.09340"
Set-TextValue "D16" "3.838"
Set-TextValue "D17" "0.001642"
Set-TextValue "D18" "0.04737"
Set-TextValue "D19" "0.006374"
Set-TextValue "D20" "0.005673"
Set-TextValue "D21" "0.001078"
Set-TextValue "D23" "3.719"
Set-TextValue "D24" "2.413"
Set-TextValue "D26" "0.1240"
Set-TextValue "D27" "0.0002704"
Set-TextValue "D40" "0.04694"
Set-TextValue "D41" "0.007011"

# ----- Rows 42 / 43 swap (CEJI <-> BKEXToken) with refreshed prices -----
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1162"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003290"
$ws.Range("E43").Value = "42CEJICEJI"

# ----- More column D updates -----
Set-TextValue "D44" "0.01157"
Set-TextValue "D45" "0.00006265"
Set-TextValue "D46" "0.0009902"

# ----- Rows 49 / 50 swap (BOLO <-> CryptobidCoin) with refreshed prices -----
$ws.Range("B49").Value = "CryptobidCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
Set-TextValue "D49" "0.00001400"
$ws.Range("E49").Value = "48CryptobidCoinCBCWorstin24h"

$ws.Range("B50").Value = "BOLO"
$ws.Range("C50").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue "D50" "0.002236"
$ws.Range("E50").Value = "49BOLOBOLO"
